$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 3.5
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 3.5
$ws.Range("R2").Value = 1.3
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05
$ws.Range("U2").Value = 2.63
$ws.Range("V2").Value = 1.44

# Row 3
$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 2
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 6.5
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 12
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 6.5
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 26
$ws.Range("AJ3").Value = 67
$ws.Range("AL3").Value = 51
$ws.Range("AO3").Value = 9.5
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 34
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AW3").Value = 7
$ws.Range("BA3").Value = 151

# Row 4
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 4.75
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 17
$ws.Range("Y4").Value = 15
$ws.Range("Z4").Value = 41
$ws.Range("AA4").Value = 41
$ws.Range("AB4").Value = 51
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 21
$ws.Range("AG4").Value = 5.5
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 10
$ws.Range("AJ4").Value = 19
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 5.5
$ws.Range("AO4").Value = 23
$ws.Range("AP4").Value = 41
$ws.Range("AQ4").Value = 81
$ws.Range("AR4").Value = 151
$ws.Range("AW4").Value = 4
$ws.Range("AX4").Value = 13
$ws.Range("AZ4").Value = 41

# Row 6
$ws.Range("G6").Value = 1.52
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 5.5
$ws.Range("J6").Value = 2
$ws.Range("L6").Value = 5.2
$ws.Range("N6").Value = 8.75
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 4.1
$ws.Range("Q6").Value = 1.6
$ws.Range("R6").Value = 2.2
$ws.Range("S6").Value = 1.31
$ws.Range("T6").Value = 3.15
$ws.Range("U6").Value = 1.7
$ws.Range("V6").Value = 2.02
$ws.Range("W6").Value = 8.25
$ws.Range("X6").Value = 7.9
$ws.Range("Y6").Value = 8
$ws.Range("Z6").Value = 11
$ws.Range("AA6").Value = 11.25
$ws.Range("AB6").Value = 22
$ws.Range("AC6").Value = 8.75
$ws.Range("AD6").Value = 8.25
$ws.Range("AE6").Value = 15.5
$ws.Range("AF6").Value = 60
$ws.Range("AG6").Value = 17.5
$ws.Range("AH6").Value = 35
$ws.Range("AI6").Value = 17
$ws.Range("AJ6").Value = 100
$ws.Range("AK6").Value = 50
$ws.Range("AL6").Value = 45
$ws.Range("AM6").Value = 400
$ws.Range("AN6").Value = 3.5
$ws.Range("AO6").Value = 7.1
$ws.Range("AP6").Value = 15
$ws.Range("AQ6").Value = 21
$ws.Range("AS6").Value = 175
$ws.Range("AT6").Value = 3.15
$ws.Range("AU6").Value = 7.4
$ws.Range("AV6").Value = 60
$ws.Range("AW6").Value = 7.2
$ws.Range("AX6").Value = 29
$ws.Range("AY6").Value = 30
$ws.Range("AZ6").Value = 175
$ws.Range("BA6").Value = 175
$ws.Range("BB6").Value = 350

# Row 7
$ws.Range("G7").Value = 2.4
$ws.Range("H7").Value = 3.05
$ws.Range("M7").Value = 1.09
$ws.Range("N7").Value = 6.2
$ws.Range("P7").Value = 2.75
$ws.Range("W7").Value = 7
$ws.Range("X7").Value = 11
$ws.Range("Z7").Value = 25
$ws.Range("AC7").Value = 6.2
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 80
$ws.Range("AH7").Value = 14.5
$ws.Range("AJ7").Value = 37
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 40
$ws.Range("AR7").Value = 90
$ws.Range("AY7").Value = 25
$ws.Range("BB7").Value = 350
